$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers I1 (I0) and J1 (IF), matching the existing header style (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-74
$ijData = @(
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(10, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(4, 4),
    @(3, 3)
)

for ($i = 0; $i -lt $ijData.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $ijData[$i][0]
    $ws.Cells.Item($r, 10).Value = $ijData[$i][1]
}
